$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B20").Value = 0.004987
$ws.Range("C20").Value = 324.21
$ws.Range("D20").Value = 360
$ws.Range("E20").Formula = "=B20*D20/C20"

$ws.Range("B22").Value = 0.004987

$ws.Range("E20").Select()
